# Refresh the LR-pair table (Bmp4 -> Bmpr1b) with newly recomputed TPM
# values. The sending/target cluster grid grows from ECs x {FAPs,MuSCs}
# (7 data rows) to ECs x {FAPs,MuSCs,Resolving-Mac} (a full 3x3 grid,
# 9 rows after the header) by introducing a new "Resolving-Mac" target
# cluster, and every row's ligand/receptor/edge statistics are updated
# to the new TPM-derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the various expression / specificity metrics.
$rows = @(
  @("ECs",   "Bmp4", "Bmpr1b", "FAPs",          3, 1, 8.133160666666667, 24.399482, 0.3870696756706061, 0.3870696756706061, 3, 1, 3.510190333333334, 10.530571, 0.7927501397588634, 0.7927501397588635, 28.54894195158023, 256.940477564222, 0.3068495394842909, 0.306849539484291),
  @("ECs",   "Bmp4", "Bmpr1b", "MuSCs",         3, 1, 8.133160666666667, 24.399482, 0.3870696756706061, 0.3870696756706061, 3, 1, 0.8869683333333332, 2.660905, 0.2003151406163121, 0.2003151406163121, 7.213855961245555, 64.92470365121, 0.07753591651026778, 0.0775359165102678),
  @("ECs",   "Bmp4", "Bmpr1b", "Resolving-Mac", 3, 1, 8.133160666666667, 24.399482, 0.3870696756706061, 0.3870696756706061, 1, 0.3333333333333333, 0.030706, 0.09211800000000001, 0.006934719624824425, 0.006934719624824427, 0.2497368314306667, 2.247631482876, 0.002684219676047378, 0.002684219676047379),
  @("FAPs",  "Bmp4", "Bmpr1b", "FAPs",          3, 1, 7.480073, 22.440219, 0.3559882250904906, 0.3559882250904906, 3, 1, 3.510190333333334, 10.530571, 0.7927501397588634, 0.7927501397588635, 26.25647993722767, 236.308319435049, 0.2822097151929961, 0.2822097151929961),
  @("FAPs",  "Bmp4", "Bmpr1b", "MuSCs",         3, 1, 7.480073, 22.440219, 0.3559882250904906, 0.3559882250904906, 3, 1, 0.8869683333333332, 2.660905, 0.2003151406163121, 0.2003151406163121, 6.634587882021666, 59.71129093819499, 0.07130983136675298, 0.07130983136675298),
  @("FAPs",  "Bmp4", "Bmpr1b", "Resolving-Mac", 3, 1, 7.480073, 22.440219, 0.3559882250904906, 0.3559882250904906, 1, 0.3333333333333333, 0.030706, 0.09211800000000001, 0.006934719624824425, 0.006934719624824427, 0.229683121538, 2.067148093842, 0.00246867853074144, 0.002468678530741441),
  @("MuSCs", "Bmp4", "Bmpr1b", "FAPs",          3, 1, 5.398902333333333, 16.196707, 0.2569420992389033, 0.2569420992389034, 3, 1, 3.510190333333334, 10.530571, 0.7927501397588634, 0.7927501397588635, 18.95117478107744, 170.560573029697, 0.2036908850815764, 0.2036908850815764),
  @("MuSCs", "Bmp4", "Bmpr1b", "MuSCs",         3, 1, 5.398902333333333, 16.196707, 0.2569420992389033, 0.2569420992389034, 3, 1, 0.8869683333333332, 2.660905, 0.2003151406163121, 0.2003151406163121, 4.788655404426111, 43.09789863983499, 0.05146939273929133, 0.05146939273929135),
  @("MuSCs", "Bmp4", "Bmpr1b", "Resolving-Mac", 3, 1, 5.398902333333333, 16.196707, 0.2569420992389033, 0.2569420992389034, 1, 0.3333333333333333, 0.030706, 0.09211800000000001, 0.006934719624824425, 0.006934719624824427, 0.1657786950473333, 1.492008255426, 0.001781821418035608, 0.001781821418035609)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $rowValues = $rows[$i]
  for ($c = 0; $c -lt $rowValues.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
  }
}
